# allocation rule updated with 5 and 10 mi rad
# Adds "Within 5 miles" and "Within 10 miles of HFC production facility" columns
# (F and G) to both the "Means" and "Standard Deviations" sheets, and refreshes
# the statistics that shifted as a result of including the new radii.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: "Means"
# ---------------------------------------------------------------------------
$means = $wb.Worksheets.Item("Means")

# New header cells for the 5-mile / 10-mile columns.
$means.Cells.Item(1, 6).Value = "Within 5 miles of HFC production facility"
$means.Cells.Item(1, 7).Value = "Within 10 miles of HFC production facility"

# Row 2: % White
$means.Cells.Item(2, 6).Value = 64
$means.Cells.Item(2, 7).Value = 74

# Row 3: % Black or African American
$means.Cells.Item(3, 6).Value = 25
$means.Cells.Item(3, 7).Value = 18

# Row 4: % Other
$means.Cells.Item(4, 6).Value = 11
$means.Cells.Item(4, 7).Value = 8.6

# Row 5: % Hispanic
$means.Cells.Item(5, 6).Value = 14
$means.Cells.Item(5, 7).Value = 9.4

# Row 6: Median Income [1,000 2019$]
$means.Cells.Item(6, 6).Value = 51
$means.Cells.Item(6, 7).Value = 58

# Row 7: % Below Poverty Line
$means.Cells.Item(7, 6).Value = 12
$means.Cells.Item(7, 7).Value = 10

# Row 8: % Below Half the Poverty Line
$means.Cells.Item(8, 6).Value = 6.4
$means.Cells.Item(8, 7).Value = 5.7

# Row 9: Total Cancer Risk (per million) -- existing columns were also
# recalculated once the 5/10 mile radii were added.
$means.Cells.Item(9, 2).Value = 29
$means.Cells.Item(9, 3).Value = 34
$means.Cells.Item(9, 4).Value = 57
$means.Cells.Item(9, 5).Value = 43
$means.Cells.Item(9, 6).Value = 39
$means.Cells.Item(9, 7).Value = 35

# Row 10: Total Respiratory (hazard quotient) -- existing columns also changed.
$means.Cells.Item(10, 2).Value = 0.37
$means.Cells.Item(10, 3).Value = 0.47
$means.Cells.Item(10, 4).Value = 0.67
$means.Cells.Item(10, 5).Value = 0.53
$means.Cells.Item(10, 6).Value = 0.48
$means.Cells.Item(10, 7).Value = 0.45

# ---------------------------------------------------------------------------
# Sheet 2: "Standard Deviations"
# ---------------------------------------------------------------------------
$sd = $wb.Worksheets.Item("Standard Deviations")

# New header cells for the 5-mile / 10-mile SD columns.
$sd.Cells.Item(1, 6).Value = "Within 5 mile of HFC production facility SD"
$sd.Cells.Item(1, 7).Value = "Within 10 mile of HFC production facility SD"

# Row 2: % White
$sd.Cells.Item(2, 6).Value = 21
$sd.Cells.Item(2, 7).Value = 21

# Row 3: % Black or African American
$sd.Cells.Item(3, 6).Value = 22
$sd.Cells.Item(3, 7).Value = 20

# Row 4: % Other
$sd.Cells.Item(4, 6).Value = 9.9
$sd.Cells.Item(4, 7).Value = 9.8

# Row 5: % Hispanic
$sd.Cells.Item(5, 6).Value = 14
$sd.Cells.Item(5, 7).Value = 13

# Row 6: Median Income [1,000 2019$]
$sd.Cells.Item(6, 6).Value = 21
$sd.Cells.Item(6, 7).Value = 22

# Row 7: % Below Poverty Line
$sd.Cells.Item(7, 6).Value = 11
$sd.Cells.Item(7, 7).Value = 11

# Row 8: % Below Half the Poverty Line
$sd.Cells.Item(8, 6).Value = 8.7
$sd.Cells.Item(8, 7).Value = 7.1

# Row 9: Total Cancer Risk (per million) -- existing columns were also
# recalculated once the 5/10 mile radii were added.
$sd.Cells.Item(9, 2).Value = 10
$sd.Cells.Item(9, 3).Value = 5.3
$sd.Cells.Item(9, 4).Value = 14
$sd.Cells.Item(9, 5).Value = 9.5
$sd.Cells.Item(9, 6).Value = 6.2
$sd.Cells.Item(9, 7).Value = 6.5

# Row 10: Total Respiratory (hazard quotient) -- existing columns also changed.
$sd.Cells.Item(10, 2).Value = 0.14
$sd.Cells.Item(10, 3).Value = 0.07
$sd.Cells.Item(10, 4).Value = 0.14
$sd.Cells.Item(10, 5).Value = 0.095
$sd.Cells.Item(10, 6).Value = 0.067
$sd.Cells.Item(10, 7).Value = 0.066
